$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Fix the SamplesTab query (cell B3): the "Tumor" column previously coalesced
# the collected/aliased `tumor` variable; it now reads the sample's
# sample_tumor_status property directly.
$newQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["LCCC 1108: Development of a Tumor Molecular Analyses Program and Its Use to Support Treatment Decisions (UNCseqTM)"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newQuery

# Move the active selection, matching the saved view state.
$ws.Range("C13").Select()
